# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets (rows captured 2026-01-28 16:53:30 - 16:54:29).

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [object]$ws,
        [int]$startRow,
        [object[]]$rows
    )

    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $rowData = $rows[$i]

        # Force the whole row to Text format first so Excel doesn't
        # reinterpret date-, time- or percent-looking strings as numbers.
        $rng = $ws.Range("A" + $r + ":F" + $r)
        $rng.NumberFormat = "@"

        for ($c = 0; $c -lt $rowData.Count; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
        }

        # Drop back to the workbook's normal (unstyled) look, matching
        # every other data row in the log.
        $rng.Style = "Normal"
    }
}

# --- PIR sheet: rows 197-210 ---------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @("2026-01-28","16:53:30","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:31","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:34","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:39","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:44","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:49","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:54","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:53:59","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:54:05","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:54:09","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:54:14","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:54:19","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:54:25","16:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","16:54:29","16:00","Bathroom","No Motion","Inactive")
)
Add-LogRows $wsPir 197 $pirRows

# --- Humidity sheet: rows 194-204 -----------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @("2026-01-28","16:53:31","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:53:39","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:53:47","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:53:51","16:00","Bathroom","87.8%","Active"),
    @("2026-01-28","16:54:03","16:00","Bathroom","86.4%","Active"),
    @("2026-01-28","16:54:07","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:54:11","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:54:15","16:00","Bathroom","87.0%","Active"),
    @("2026-01-28","16:54:20","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:54:23","16:00","Bathroom","87.9%","Active"),
    @("2026-01-28","16:54:27","16:00","Bathroom","87.0%","Active")
)
Add-LogRows $wsHumidity 194 $humidityRows

# --- Temperature sheet: rows 194-204 ---------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @("2026-01-28","16:53:32","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:53:40","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:53:48","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:53:52","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:04","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:08","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:12","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:16","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:20","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:24","16:00","Bathroom","22.8C","Active"),
    @("2026-01-28","16:54:28","16:00","Bathroom","22.8C","Active")
)
Add-LogRows $wsTemperature 194 $temperatureRows
